$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.880.70"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.859.03"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'304.20"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "'0.3627"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("D9").Value = "'0.07160"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.8914"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.856.36"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "'0.07453"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "'92.86"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").Value = "'5.225"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.000008493"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "'14.07"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "26.921.05"
$ws.Range("D21").Value = "'5.013"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").Value = "2.083.69"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").Value = "'6.417"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").Value = "'147.55"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value = "'1.793"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "'2.045"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").Value = "'113.10"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'4.649"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").Value = "'4.658"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'0.09238"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").Value = "'0.05076"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").Value = "'0.7447"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "'2.972"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'3.263"
$ws.Range("E37").Value = "  +6.99%  "
$ws.Range("D38").Value = "'2.514"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "'1.086"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").Value = "'0.5336"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").Value = "'118.46"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").Value = "'6.495"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "'8.438"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").Value = "'0.1460"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "'0.4635"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "'0.9994"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "'10.03"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'36.77"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'62.86"
$ws.Range("E51").Value = "  -3.27%  "
